# Consolidate the time_counts sheet: merge related area labels (e.g. K2/K1/K -> K,
# J3/J2/J1/J -> J, Pg/K2-Pg -> Pg, Mz/J1-J2 -> Mz) and drop the now-redundant rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: K2 (1271) -> K (2522)
$ws.Range("A2").Value = "K"
$ws.Range("B2").Value = 2522

# Row 3: K1 (1152) -> J (825)
$ws.Range("A3").Value = "J"
$ws.Range("B3").Value = 825

# Row 4: Tr (517) stays as-is

# Row 5: J3 (411) -> Pg (365)
$ws.Range("A5").Value = "Pg"
$ws.Range("B5").Value = 365

# Row 6: Tr-J (338) stays as-is

# Row 7: Pg (194) -> Mz (88)
$ws.Range("A7").Value = "Mz"
$ws.Range("B7").Value = 88

# The data that used to live in rows 8-16 has been folded into the rows above,
# so remove those now-duplicate rows entirely (shifts dimension down to A1:B7).
$ws.Rows("8:16").Delete()
